$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A id, B game_class_id, C name, D description, E requires_class_rank_level,
# F specialty_damage, G increase_specialty_damage_per_level,
# H specialty_damage_uses_damage_stat_amount, I base_damage_mod, J base_ac_mod,
# K base_healing_mod, L base_spell_damage_mod, M health_mod,
# N base_damage_stat_increase, O attack_type_required

# Row 136 (id 135)
$ws.Cells.Item(136, 1).Value = 135
$ws.Cells.Item(136, 2).Value = 16
$ws.Cells.Item(136, 3).Value = "Pale Whisper"
$ws.Cells.Item(136, 4).Value = "Increase your Damage and Healing by +50% at level 100."
$ws.Cells.Item(136, 5).Value = 1
$ws.Cells.Item(136, 9).Value = 0.05
$ws.Cells.Item(136, 11).Value = 0.05

# Row 137 (id 136)
$ws.Cells.Item(137, 1).Value = 136
$ws.Cells.Item(137, 2).Value = 16
$ws.Cells.Item(137, 3).Value = "Tincture of Faith"
$ws.Cells.Item(137, 4).Value = "Grow your healing by +200% at level 100."
$ws.Cells.Item(137, 5).Value = 12
$ws.Cells.Item(137, 11).Value = 0.02

# Row 138 (id 137)
$ws.Cells.Item(138, 1).Value = 137
$ws.Cells.Item(138, 2).Value = 16
$ws.Cells.Item(138, 3).Value = "Hollow Existance"
$ws.Cells.Item(138, 4).Value = "Grow your damage by +75%, your AC (defence) by +100% and your healing by +250% at level 100."
$ws.Cells.Item(138, 5).Value = 24
$ws.Cells.Item(138, 9).Value = 0.0075
$ws.Cells.Item(138, 10).Value = 0.01
$ws.Cells.Item(138, 11).Value = 0.025

# Row 139 (id 138)
$ws.Cells.Item(139, 1).Value = 138
$ws.Cells.Item(139, 2).Value = 16
$ws.Cells.Item(139, 3).Value = "Chalice of suffering"
$ws.Cells.Item(139, 4).Value = "Increase damage by +100%, Healing by +250% and Health by +200% at level 100"
$ws.Cells.Item(139, 5).Value = 36
$ws.Cells.Item(139, 11).Value = 0.025
$ws.Cells.Item(139, 13).Value = 0.02
$ws.Cells.Item(139, 14).Value = 0.01

# Row 140 (id 139)
$ws.Cells.Item(140, 1).Value = 139
$ws.Cells.Item(140, 2).Value = 16
$ws.Cells.Item(140, 3).Value = "Cornered in a fight"
$ws.Cells.Item(140, 4).Value = "Increase your damage by +300% at level 100"
$ws.Cells.Item(140, 5).Value = 48
$ws.Cells.Item(140, 14).Value = 0.03

# Row 141 (id 140)
$ws.Cells.Item(141, 1).Value = 140
$ws.Cells.Item(141, 2).Value = 16
$ws.Cells.Item(141, 3).Value = "Faithless Aboration"
$ws.Cells.Item(141, 4).Value = "Increase your healing by +300% at level 100 and your spell damage by +150% at level 100."
$ws.Cells.Item(141, 5).Value = 60
$ws.Cells.Item(141, 11).Value = 0.03
$ws.Cells.Item(141, 12).Value = 0.015

# Row 142 (id 141)
$ws.Cells.Item(142, 1).Value = 141
$ws.Cells.Item(142, 2).Value = 16
$ws.Cells.Item(142, 3).Value = "Graverobbing Shadows of Pain"
$ws.Cells.Item(142, 4).Value = "Deal damage equal to 2500 + 2.5% of your damage, growing by a total of +2500 at level 100"
$ws.Cells.Item(142, 5).Value = 70
$ws.Cells.Item(142, 6).Value = 2500
$ws.Cells.Item(142, 7).Value = 25
$ws.Cells.Item(142, 8).Value = 0.025
$ws.Cells.Item(142, 15).Value = "attack"

# Row 143 (id 142)
$ws.Cells.Item(143, 1).Value = 142
$ws.Cells.Item(143, 2).Value = 16
$ws.Cells.Item(143, 3).Value = "Necroctic Emotional Drift"
$ws.Cells.Item(143, 4).Value = "Deal 8,000 damage + 6% of your damage stat as damage with a bonus of +8,000 damage at level 100. You will also grow your damage modifier by +200% at level 100."
$ws.Cells.Item(143, 5).Value = 80
$ws.Cells.Item(143, 6).Value = 8000
$ws.Cells.Item(143, 7).Value = 80
$ws.Cells.Item(143, 8).Value = 0.06
$ws.Cells.Item(143, 9).Value = 0.02
$ws.Cells.Item(143, 15).Value = "attack_and_cast"

# Row 144 (id 143)
$ws.Cells.Item(144, 1).Value = 143
$ws.Cells.Item(144, 2).Value = 16
$ws.Cells.Item(144, 3).Value = "The line between life and death"
$ws.Cells.Item(144, 4).Value = "Deal 20,000 damage while use 20% of yuor damage stat, growing to +20,000 damage at level 100. You will also grow your damage stat by +200% at level 100."
$ws.Cells.Item(144, 5).Value = 90
$ws.Cells.Item(144, 6).Value = 20000
$ws.Cells.Item(144, 7).Value = 200
$ws.Cells.Item(144, 8).Value = 0.12
$ws.Cells.Item(144, 9).Value = 0.02
$ws.Cells.Item(144, 14).Value = 0.02
$ws.Cells.Item(144, 15).Value = "cast_and_attack"
